$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -5
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = -3
$ws.Range("F15").Value = -3
$ws.Range("F20").Value = 1
$ws.Range("F23").Value = -1
$ws.Range("F25").Value = -1
$ws.Range("F29").Value = -4
$ws.Range("F30").Value = 0
$ws.Range("F32").Value = -9
$ws.Range("F34").Value = 2
$ws.Range("F35").Value = 2
$ws.Range("F37").Value = -8
$ws.Range("F39").Value = -14
$ws.Range("F43").Value = -1
$ws.Range("F45").Value = -4
